$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: RowNumber, Coin, Link, Price, Volume(1h)
$data = @(
    @(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "29.197.87", "  +0.45%  "),
    @(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "1.831.08", "  -0.20%  "),
    @(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "0.9986", "  -0.12%  "),
    @(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "243.04", "  +0.19%  "),
    @(6, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.6210", "  +1.30%  "),
    @(7, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.001", "  +0.00%  "),
    @(8, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.07357", "  -1.36%  "),
    @(9, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.2912", "  -0.22%  "),
    @(10, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "23.22", "  +0.48%  "),
    @(11, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.07665", "  -0.18%  "),
    @(12, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "1.837.05", "  +0.03%  "),
    @(13, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "4.973", "  -0.66%  "),
    @(14, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "0.6689", "  -0.32%  "),
    @(15, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "82.54", "  -0.05%  "),
    @(16, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.000009020", "  -1.15%  "),
    @(17, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "5.848", "  -1.28%  "),
    @(18, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "29.178.76", "  +0.41%  "),
    @(19, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "2.088.33", "  -0.05%  "),
    @(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "235.75", "  +2.09%  "),
    @(21, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "12.48", "  -1.33%  "),
    @(22, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "1.001", "  -0.04%  "),
    @(23, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "7.370", "  +2.57%  "),
    @(24, "BinanceUSD", "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd", "1.002", "  +0.09%  "),
    @(25, "Monero", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr", "158.40", "  -0.68%  "),
    @(26, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "8.550", "  +0.71%  "),
    @(27, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.1392", "  +0.49%  "),
    @(28, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "17.63", "  -0.85%  "),
    @(29, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "1.493", "  -0.23%  "),
    @(30, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.05823", "  +4.85%  "),
    @(31, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "4.085", "  -1.19%  "),
    @(32, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "4.090", "  -1.51%  "),
    @(33, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "1.215", "  +0.83%  "),
    @(34, "LidoDAOToken", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo", "1.868", "  +1.67%  "),
    @(35, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "0.7272", "  -2.20%  "),
    @(36, "ARBITRUM", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb", "1.139", "  -0.04%  "),
    @(37, "HuobiToken", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht", "2.608", "  -1.95%  "),
    @(38, "MXToken", "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx", "2.852", "  +2.94%  "),
    @(39, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "1.222.71", "  +1.30%  "),
    @(40, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.01752", "  -1.56%  "),
    @(41, "FraxShare", "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs", "6.244", "  -3.68%  "),
    @(42, "TrustWalletToken", "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt", "0.9096", "  +2.17%  "),
    @(43, "PaxDollar", "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp", "1.002", "  +0.15%  "),
    @(44, "RocketPoolETH", "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth", "1.991.69", "  +0.46%  "),
    @(45, "Quant", "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt", "101.74", "  -0.21%  "),
    @(46, "Aave", "https://coinranking.com/coin/ixgUfzmLR+aave-aave", "65.50", "  -0.08%  "),
    @(47, "BabyDogeCoin", "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge", "0.00000000120", "  -1.86%  "),
    @(48, "Mantle", "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt", "0.5046", "  -0.83%  "),
    @(49, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.4029", "  -0.84%  "),
    @(50, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "9.122", "  +0.19%  "),
    @(51, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1131", "  +2.88%  ")
)

foreach ($item in $data) {
    $r = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    # Force Price and Volume columns to text so numeric-looking
    # strings (e.g. "1.001", "29.197.87") aren't coerced to numbers.
    $ws.Cells.Item($r, 4).NumberFormat = "@"
    $ws.Cells.Item($r, 4).Value = $item[3]
    $ws.Cells.Item($r, 5).NumberFormat = "@"
    $ws.Cells.Item($r, 5).Value = $item[4]
}
